$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: update existing values and add new ones (F2:I2 previously empty)
$ws.Range("B2").Value = 0.9
$ws.Range("C2").Value = 0.58
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0.9
$ws.Range("F2").Value = 1.62
$ws.Range("G2").Value = 3.01
$ws.Range("H2").Value = 0.9
$ws.Range("I2").Value = 0.38

# Row 3: update D3 and E3
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.9

# Row 4: update D4 and E4
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0.9
